$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Tabel1" story backlog (H1:M28, plus an inline N note) gets
# reshuffled: three stories ("thema", "limiet van storten " and
# "break message") are marked done and bubble up near the top of the
# list, while the rest shift down one row, keeping their relative
# order. Row 21 ends up holding exactly what row 20 used to hold
# (including its "?" effective-duration cell and its N-column note),
# so write row 21's values first and copy row 20's current (still
# untouched) formatting onto it before row 20 itself gets overwritten
# with its own new content.
# ------------------------------------------------------------------

# Row 21: "Mogelijke winst zien" (was row 20, still open, keeps the
# "?" effective-duration value and the N-column note + formatting)
$ws.Range("H21").Value2 = 20
$ws.Range("I21").Value2 = "Mogelijke winst zien"
$ws.Range("J21").Value2 = 0.5
$ws.Range("K21").Value2 = "?"
$ws.Range("L21").Value2 = "Rein"
$ws.Range("M21").Value2 = "Ja"
$ws.Range("N21").Value2 = "--> niet aangevuld "

$ws.Range("K20").Copy()
$ws.Range("K21").PasteSpecial(-4122)
$ws.Range("N20").Copy()
$ws.Range("N21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 19: "thema" -> done (Sepp, 1.5/1.5, "ja")
$ws.Range("H19").Value2 = 18
$ws.Range("I19").Value2 = "thema"
$ws.Range("J19").Value2 = 1.5
$ws.Range("K19").Value2 = 1.5
$ws.Range("L19").Value2 = "Sepp"
$ws.Range("M19").Value2 = "ja"

# Row 20: "limiet van storten " -> done (Sepp, 0.5/0.25, "Ja")
# (row 20 used to be the "?" / red-font row before the reshuffle, so
# its old K-column formatting must be dropped along with the N note)
$ws.Range("H20").Value2 = 19
$ws.Range("I20").Value2 = "limiet van storten "
$ws.Range("J20").Value2 = 0.5
$ws.Range("K20").Value2 = 0.25
$ws.Range("K20").ClearFormats()
$ws.Range("L20").Value2 = "Sepp"
$ws.Range("M20").Value2 = "Ja"
$ws.Range("N20").Clear()

# Row 22: "break message" -> done (Sepp, 0.25/0.25, "ja")
$ws.Range("H22").Value2 = 21
$ws.Range("I22").Value2 = "break message"
$ws.Range("J22").Value2 = 0.25
$ws.Range("K22").Value2 = 0.25
$ws.Range("L22").Value2 = "Sepp"
$ws.Range("M22").Value2 = "ja"

# Rows 23-27: remaining backlog items shift down one row, still open
$ws.Range("H23").Value2 = 22
$ws.Range("I23").Value2 = "geld verloren lijsten"
$ws.Range("J23").Value2 = 0.25
$ws.Range("K23").Value2 = 0
$ws.Range("M23").Value2 = "Nee"

$ws.Range("H24").Value2 = 23
$ws.Range("I24").Value2 = "overzicht vorige races"
$ws.Range("J24").Value2 = 2
$ws.Range("K24").Value2 = 0
$ws.Range("M24").Value2 = "Nee"

$ws.Range("H25").Value2 = 24
$ws.Range("I25").Value2 = "paarden beheren "
$ws.Range("J25").Value2 = 4
$ws.Range("K25").Value2 = 0
$ws.Range("M25").Value2 = "Nee"

$ws.Range("H26").Value2 = 25
$ws.Range("I26").Value2 = "track veranderen(functioneel)"
$ws.Range("J26").Value2 = 2
$ws.Range("K26").Value2 = 0
$ws.Range("M26").Value2 = "Nee"

$ws.Range("H27").Value2 = 26
$ws.Range("I27").Value2 = "grafiek overzicht"
$ws.Range("J27").Value2 = 6
$ws.Range("K27").Value2 = 0
$ws.Range("M27").Value2 = "Nee"

# Selection moves to L22, matching the edited workbook's cursor position
$ws.Range("L22").Select()

$wb.Save()
